{"js": "// Replace the \"RPC Explorer\" menu item text with \"Insight Explorer\".\nconst body = context.document.body;\n\nconst searchResults = body.search(\"RPC Explorer\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"Insight Explorer\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the \"RPC Explorer\" menu item text with \"Insight Explorer\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"RPC Explorer\"\n$find.Replacement.Text = \"Insight Explorer\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n\n# wdFindContinue=1, wdReplaceAll=2\n$find.Execute(\"RPC Explorer\", $false, $false, $false, $false, $false, $true, 1, $false, \"Insight Explorer\", 2) | Out-Null\n"}
